$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 97) to the daily log sheet, mirroring the
# existing rows: date (as text), weekday (as text), hour, ranking.
#
# Column A holds date-looking strings ("2025/10/13") but they are stored
# as plain text in this workbook (no number format), so we force a text
# number format before assignment to stop Excel from auto-converting the
# literal into a real date serial, then clear the format again afterward
# so the new cell carries no explicit style, matching the rest of the
# sheet's unstyled data cells.
$ws.Range("A97").NumberFormat = "@"
$ws.Range("A97").Value = "2025/10/13"
$ws.Range("A97").ClearFormats()

$ws.Range("B97").Value = "月"
$ws.Range("C97").Value = 5
$ws.Range("D97").Value = 201
